$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 58
$ws.Range("AI58").Value = 556231.04
$ws.Range("AJ58").Value = 292852
$ws.Range("AK58").Value = 220658
$ws.Range("AL58").Value = 209483.008
$ws.Range("AM58").Value = 396392.96
$ws.Range("AN58").Value = 238403.008
$ws.Range("AO58").Value = 230744
$ws.Range("AP58").Value = 272468.992
$ws.Range("AQ58").Value = 237273.024
$ws.Range("AR58").Value = 214184.992
$ws.Range("CH58").Value = 419139.008
$ws.Range("CI58").Value = 487592.928
$ws.Range("CJ58").Value = 275516
$ws.Range("CK58").Value = 892444.032
$ws.Range("CL58").Value = 791254.976
$ws.Range("CM58").Value = 751369.088
$ws.Range("CN58").Value = 733523.968
$ws.Range("CO58").Value = 729696
$ws.Range("CP58").Value = 895382.0159999999
$ws.Range("CQ58").Value = 719520

# Row 59
$ws.Range("AI59").Value = -452532.032
$ws.Range("AJ59").Value = -262560
$ws.Range("AK59").Value = -189426
$ws.Range("AL59").Value = -173792.992
$ws.Range("AM59").Value = -283889.984
$ws.Range("AN59").Value = -189932
$ws.Range("AO59").Value = -186180
$ws.Range("AP59").Value = -216520.992
$ws.Range("AQ59").Value = -176264.032
$ws.Range("AR59").Value = -188772
$ws.Range("CH59").Value = -377734.016
$ws.Range("CI59").Value = -458587.968
$ws.Range("CJ59").Value = -235900
$ws.Range("CK59").Value = -844428.032
$ws.Range("CL59").Value = -711283.968
$ws.Range("CM59").Value = -747913.856
$ws.Range("CN59").Value = -728374.0159999999
$ws.Range("CO59").Value = -686873.9840000001
$ws.Range("CP59").Value = -868883.008
$ws.Range("CQ59").Value = -788768.064

# Row 60
$ws.Range("AI60").Value = 103699
$ws.Range("AJ60").Value = 30292
$ws.Range("AK60").Value = 31232
$ws.Range("AL60").Value = 35690
$ws.Range("AM60").Value = 112502.992
$ws.Range("AN60").Value = 48471
$ws.Range("AO60").Value = 44564
$ws.Range("AP60").Value = 55948
$ws.Range("AQ60").Value = 61009
$ws.Range("AR60").Value = 25413
$ws.Range("CH60").Value = 41405
$ws.Range("CI60").Value = 29005.008
$ws.Range("CJ60").Value = 39616
$ws.Range("CK60").Value = 48016
$ws.Range("CL60").Value = 79971
$ws.Range("CM60").Value = 3455
$ws.Range("CN60").Value = 5150
$ws.Range("CO60").Value = 42822
$ws.Range("CP60").Value = 26499
$ws.Range("CQ60").Value = -69248

# Row 61
$ws.Range("AI61").Value = 2472
$ws.Range("AJ61").Value = -5053
$ws.Range("AK61").Value = -12913
$ws.Range("AL61").Value = -4549
$ws.Range("AM61").Value = -81850
$ws.Range("AN61").Value = -6115
$ws.Range("AO61").Value = -14707
$ws.Range("AP61").Value = -28666
$ws.Range("AQ61").Value = -49239
$ws.Range("AR61").Value = -7173
$ws.Range("CH61").ClearContents()
$ws.Range("CI61").ClearContents()
$ws.Range("CJ61").ClearContents()
$ws.Range("CK61").ClearContents()
$ws.Range("CL61").ClearContents()
$ws.Range("CM61").ClearContents()
$ws.Range("CN61").ClearContents()
$ws.Range("CO61").ClearContents()
$ws.Range("CP61").ClearContents()
$ws.Range("CQ61").ClearContents()

# Row 62
$ws.Range("AI62").Value = 25341
$ws.Range("AJ62").Value = 4820
$ws.Range("AK62").Value = 5773
$ws.Range("AL62").Value = 6647
$ws.Range("AM62").Value = 18235
$ws.Range("AN62").Value = 7857
$ws.Range("AO62").Value = 7594
$ws.Range("AP62").Value = 7374
$ws.Range("AQ62").Value = 4689
$ws.Range("AR62").Value = 5969
$ws.Range("CH62").Value = 15973
$ws.Range("CI62").Value = 25875
$ws.Range("CJ62").Value = 15175
$ws.Range("CK62").Value = 17488
$ws.Range("CL62").Value = 20247
$ws.Range("CM62").Value = 19899
$ws.Range("CN62").Value = 16924
$ws.Range("CO62").Value = 19404
$ws.Range("CP62").Value = 21399
$ws.Range("CQ62").Value = 16698

# Row 63
$ws.Range("AI63").Value = -26278
$ws.Range("AJ63").Value = -10683
$ws.Range("AK63").Value = -11032
$ws.Range("AL63").Value = -10920
$ws.Range("AM63").Value = -28119
$ws.Range("AN63").Value = -15177
$ws.Range("AO63").Value = -16297
$ws.Range("AP63").Value = -15172
$ws.Range("AQ63").Value = -21232
$ws.Range("AR63").Value = -12671
$ws.Range("CH63").Value = -28381
$ws.Range("CI63").Value = -26243
$ws.Range("CJ63").Value = -28845
$ws.Range("CK63").Value = -26190
$ws.Range("CL63").Value = -30473
$ws.Range("CM63").Value = -31332
$ws.Range("CN63").Value = -34799
$ws.Range("CO63").Value = -34020
$ws.Range("CP63").Value = -33187
$ws.Range("CQ63").Value = -35505.008

# Row 64
$ws.Range("AI64").Value = -39155
$ws.Range("AJ64").Value = -6743
$ws.Range("AK64").Value = -6886
$ws.Range("AL64").Value = -8375
$ws.Range("AM64").Value = -28107
$ws.Range("AN64").Value = -13910
$ws.Range("AO64").Value = -11252
$ws.Range("AP64").Value = -15505
$ws.Range("AQ64").Value = -16555
$ws.Range("AR64").Value = -7733
$ws.Range("CH64").Value = -10629
$ws.Range("CI64").Value = -13718
$ws.Range("CJ64").Value = -12834
$ws.Range("CK64").Value = -13494
$ws.Range("CL64").Value = -13220
$ws.Range("CM64").Value = -14017
$ws.Range("CN64").Value = -12390
$ws.Range("CO64").Value = -13088
$ws.Range("CP64").Value = -15426
$ws.Range("CQ64").Value = -13583

# Row 65
$ws.Range("AI65").Value = -5046
$ws.Range("AJ65").Value = -267
$ws.Range("AK65").Value = -287
$ws.Range("AL65").Value = -338
$ws.Range("AM65").Value = -5401
$ws.Range("AN65").Value = -1837
$ws.Range("AO65").Value = -2459
$ws.Range("AP65").Value = -1786
$ws.Range("AQ65").Value = 6082
$ws.Range("AR65").Value = -337
$ws.Range("CH65").Value = -4021
$ws.Range("CI65").Value = -4537
$ws.Range("CJ65").Value = -2799
$ws.Range("CK65").Value = -4746
$ws.Range("CL65").Value = -5092
$ws.Range("CM65").Value = -5903
$ws.Range("CN65").Value = -4660
$ws.Range("CO65").Value = -3184
$ws.Range("CP65").Value = -4012
$ws.Range("CQ65").Value = -4310

# Row 66
$ws.Range("AI66").Value = 113957
$ws.Range("AJ66").Value = 632
$ws.Range("AK66").Value = 623
$ws.Range("AL66").Value = 926
$ws.Range("AM66").Value = 33617
$ws.Range("AN66").Value = 36448
$ws.Range("AO66").Value = 17396
$ws.Range("AP66").Value = 3692
$ws.Range("AQ66").Value = 14584
$ws.Range("AR66").Value = 3354
$ws.Range("CH66").Value = 1848
$ws.Range("CI66").Value = -1631
$ws.Range("CJ66").Value = 1200
$ws.Range("CK66").Value = 722
$ws.Range("CL66").Value = 694
$ws.Range("CM66").Value = 1243
$ws.Range("CN66").Value = 3410
$ws.Range("CO66").Value = 696
$ws.Range("CP66").Value = 1093
$ws.Range("CQ66").Value = 1390

# Row 67
$ws.Range("AI67").Value = -12961
$ws.Range("AJ67").Value = -3627
$ws.Range("AK67").Value = -8625
$ws.Range("AL67").Value = -3728
$ws.Range("AM67").Value = -42500
$ws.Range("AN67").Value = -19496
$ws.Range("AO67").Value = -9689
$ws.Range("AP67").Value = -7269
$ws.Range("AQ67").Value = -36807
$ws.Range("AR67").Value = -4682
$ws.Range("CH67").Value = -7665
$ws.Range("CI67").Value = -6123
$ws.Range("CJ67").Value = -6334
$ws.Range("CK67").Value = -3488
$ws.Range("CL67").Value = -902
$ws.Range("CM67").Value = -1172
$ws.Range("CN67").Value = -785
$ws.Range("CO67").Value = -3988
$ws.Range("CP67").Value = -2007
$ws.Range("CQ67").Value = -1271

# Row 68
$ws.Range("AI68").Value = -53386
$ws.Range("AJ68").Value = 10815
$ws.Range("AK68").Value = 7521
$ws.Range("AL68").Value = 11239
$ws.Range("AM68").Value = -29575
$ws.Range("AN68").Value = 0
$ws.Range("AO68").Value = 0
$ws.Range("AP68").Value = 0
$ws.Range("AQ68").Value = 0
$ws.Range("AR68").Value = 8927
$ws.Range("CH68").Value = 9221
$ws.Range("CI68").Value = 16229
$ws.Range("CJ68").Value = 14366
$ws.Range("CK68").Value = 23912
$ws.Range("CL68").Value = 22302
$ws.Range("CM68").Value = 31087
$ws.Range("CN68").Value = 19821
$ws.Range("CO68").Value = 23595
$ws.Range("CP68").Value = 23413
$ws.Range("CQ68").Value = 26776

# Row 69
$ws.Range("AI69").Value = 106171
$ws.Range("AJ69").Value = 25239
$ws.Range("AK69").Value = 18319
$ws.Range("AL69").Value = 31141
$ws.Range("AM69").Value = 30653
$ws.Range("AN69").Value = 42356
$ws.Range("AO69").Value = 29857
$ws.Range("AP69").Value = 27282
$ws.Range("AQ69").Value = 11770
$ws.Range("AR69").Value = 18240
$ws.Range("CH69").ClearContents()
$ws.Range("CI69").ClearContents()
$ws.Range("CJ69").ClearContents()
$ws.Range("CK69").ClearContents()
$ws.Range("CL69").ClearContents()
$ws.Range("CM69").ClearContents()
$ws.Range("CN69").ClearContents()
$ws.Range("CO69").ClearContents()
$ws.Range("CP69").ClearContents()
$ws.Range("CQ69").ClearContents()

# Row 70
$ws.Range("AI70").Value = -54
$ws.Range("AJ70").Value = 139
$ws.Range("AK70").Value = 1290
$ws.Range("AL70").Value = 700
$ws.Range("AM70").Value = -2129
$ws.Range("AN70").Value = 0
$ws.Range("AO70").Value = 0
$ws.Range("AP70").Value = 0
$ws.Range("AQ70").Value = -34416
$ws.Range("AR70").Value = 112
$ws.Range("CH70").ClearContents()
$ws.Range("CI70").ClearContents()
$ws.Range("CJ70").ClearContents()
$ws.Range("CK70").ClearContents()
$ws.Range("CL70").ClearContents()
$ws.Range("CM70").ClearContents()
$ws.Range("CN70").ClearContents()
$ws.Range("CO70").ClearContents()
$ws.Range("CP70").ClearContents()
$ws.Range("CQ70").ClearContents()

# Row 71
$ws.Range("AI71").Value = -54
$ws.Range("AJ71").Value = 139
$ws.Range("AK71").Value = 1290
$ws.Range("AL71").Value = 700
$ws.Range("AM71").Value = -2129
$ws.Range("AN71").Value = 0
$ws.Range("AO71").Value = 0
$ws.Range("AP71").Value = 0
$ws.Range("AQ71").Value = -54476
$ws.Range("AR71").Value = 118
$ws.Range("CH71").ClearContents()
$ws.Range("CI71").ClearContents()
$ws.Range("CJ71").ClearContents()
$ws.Range("CK71").ClearContents()
$ws.Range("CL71").ClearContents()
$ws.Range("CM71").ClearContents()
$ws.Range("CN71").ClearContents()
$ws.Range("CO71").ClearContents()
$ws.Range("CP71").ClearContents()
$ws.Range("CQ71").ClearContents()

# Row 72
$ws.Range("AI72").Value = 0
$ws.Range("AJ72").Value = 0
$ws.Range("AK72").Value = 0
$ws.Range("AL72").Value = 0
$ws.Range("AM72").Value = 0
$ws.Range("AN72").Value = 0
$ws.Range("AO72").Value = 0
$ws.Range("AP72").Value = 0
$ws.Range("AQ72").Value = 20060
$ws.Range("AR72").Value = -6
$ws.Range("CH72").ClearContents()
$ws.Range("CI72").ClearContents()
$ws.Range("CJ72").ClearContents()
$ws.Range("CK72").ClearContents()
$ws.Range("CL72").ClearContents()
$ws.Range("CM72").ClearContents()
$ws.Range("CN72").ClearContents()
$ws.Range("CO72").ClearContents()
$ws.Range("CP72").ClearContents()
$ws.Range("CQ72").ClearContents()

# Row 73
$ws.Range("AI73").Value = 106117
$ws.Range("AJ73").Value = 25378
$ws.Range("AK73").Value = 19609
$ws.Range("AL73").Value = 31841
$ws.Range("AM73").Value = 28524
$ws.Range("AN73").Value = 42356
$ws.Range("AO73").Value = 29857
$ws.Range("AP73").Value = 27282
$ws.Range("AQ73").Value = -22646
$ws.Range("AR73").Value = 18352
$ws.Range("CH73").Value = 17751
$ws.Range("CI73").Value = 18813
$ws.Range("CJ73").Value = 19545
$ws.Range("CK73").Value = 42220
$ws.Range("CL73").Value = 73527
$ws.Range("CM73").Value = 3260
$ws.Range("CN73").Value = -7329
$ws.Range("CO73").Value = 32237
$ws.Range("CP73").Value = 17772
$ws.Range("CQ73").Value = -79053

# Row 74
$ws.Range("AI74").Value = -64007
$ws.Range("AJ74").Value = 3340
$ws.Range("AK74").Value = 1889
$ws.Range("AL74").Value = -9078
$ws.Range("AM74").Value = -19053
$ws.Range("AN74").Value = -22300
$ws.Range("AO74").Value = -10530
$ws.Range("AP74").Value = -8645
$ws.Range("AQ74").Value = 41475
$ws.Range("AR74").Value = -3405
$ws.Range("CH74").Value = -3700
$ws.Range("CI74").Value = -1568
$ws.Range("CJ74").Value = -1147
$ws.Range("CK74").Value = 1211
$ws.Range("CL74").Value = -18742
$ws.Range("CM74").Value = 24410
$ws.Range("CN74").Value = 14202
$ws.Range("CO74").Value = 3974
$ws.Range("CP74").Value = 5388
$ws.Range("CQ74").Value = 55980

# Row 75
$ws.Range("AI75").Value = 11700
$ws.Range("AJ75").Value = -7869
$ws.Range("AK75").Value = 1390
$ws.Range("AL75").Value = 1185
$ws.Range("AM75").Value = 10932
$ws.Range("AN75").Value = 7656
$ws.Range("AO75").Value = 7235
$ws.Range("AP75").Value = 3897
$ws.Range("AQ75").Value = 58061
$ws.Range("AR75").Value = -189
$ws.Range("CH75").ClearContents()
$ws.Range("CI75").ClearContents()
$ws.Range("CJ75").ClearContents()
$ws.Range("CK75").ClearContents()
$ws.Range("CL75").ClearContents()
$ws.Range("CM75").ClearContents()
$ws.Range("CN75").ClearContents()
$ws.Range("CO75").ClearContents()
$ws.Range("CP75").ClearContents()
$ws.Range("CQ75").ClearContents()

# Row 76
$ws.Range("AI76").Value = -854
$ws.Range("AJ76").Value = -1253
$ws.Range("AK76").Value = -795
$ws.Range("AL76").Value = -881
$ws.Range("AM76").Value = -2354
$ws.Range("AN76").Value = -1478
$ws.Range("AO76").Value = -1636
$ws.Range("AP76").Value = -1447
$ws.Range("AQ76").Value = -710
$ws.Range("AR76").Value = -1729
$ws.Range("CH76").Value = -1865
$ws.Range("CI76").Value = -2166
$ws.Range("CJ76").Value = -4561
$ws.Range("CK76").Value = -7887
$ws.Range("CL76").Value = -6371
$ws.Range("CM76").Value = -6825
$ws.Range("CN76").Value = -4663
$ws.Range("CO76").Value = -7431
$ws.Range("CP76").Value = -7810
$ws.Range("CQ76").Value = -7412

# Row 77
$ws.Range("AI77").Value = 0
$ws.Range("AJ77").Value = 0
$ws.Range("AK77").Value = 0
$ws.Range("AL77").Value = 0
$ws.Range("AM77").Value = 0
$ws.Range("AN77").Value = 0
$ws.Range("AO77").Value = 0
$ws.Range("AP77").Value = 0
$ws.Range("AQ77").Value = 0
$ws.Range("AR77").Value = 0
$ws.Range("CH77").ClearContents()
$ws.Range("CI77").ClearContents()
$ws.Range("CJ77").ClearContents()
$ws.Range("CK77").ClearContents()
$ws.Range("CL77").ClearContents()
$ws.Range("CM77").ClearContents()
$ws.Range("CN77").ClearContents()
$ws.Range("CO77").ClearContents()
$ws.Range("CP77").ClearContents()
$ws.Range("CQ77").ClearContents()

# Row 78
$ws.Range("AI78").ClearContents()
$ws.Range("AJ78").ClearContents()
$ws.Range("AK78").ClearContents()
$ws.Range("AL78").ClearContents()
$ws.Range("AM78").ClearContents()
$ws.Range("AN78").Value = -1995
$ws.Range("AO78").Value = -2146
$ws.Range("AP78").Value = -2804
$ws.Range("AQ78").ClearContents()
$ws.Range("AR78").ClearContents()
$ws.Range("CH78").ClearContents()
$ws.Range("CI78").ClearContents()
$ws.Range("CJ78").ClearContents()
$ws.Range("CK78").ClearContents()
$ws.Range("CL78").ClearContents()
$ws.Range("CM78").ClearContents()
$ws.Range("CN78").ClearContents()
$ws.Range("CO78").ClearContents()
$ws.Range("CP78").ClearContents()
$ws.Range("CQ78").ClearContents()

# Row 79
$ws.Range("AI79").Value = 23331
$ws.Range("AJ79").Value = 19596
$ws.Range("AK79").Value = 22093
$ws.Range("AL79").Value = 23067
$ws.Range("AM79").Value = 12027
$ws.Range("AN79").Value = 24239
$ws.Range("AO79").Value = 22780
$ws.Range("AP79").Value = 18283
$ws.Range("AQ79").Value = 16776
$ws.Range("AR79").Value = 13029
$ws.Range("CH79").Value = 12186
$ws.Range("CI79").Value = 21128
$ws.Range("CJ79").Value = 13837
$ws.Range("CK79").Value = 35544
$ws.Range("CL79").Value = 48414
$ws.Range("CM79").Value = 20845
$ws.Range("CN79").Value = 2210
$ws.Range("CO79").Value = 28780
$ws.Range("CP79").Value = 15350
$ws.Range("CQ79").Value = -30485
